$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 5).Value = "extra_info: {""date"": {""comment"": ""ou 15600325 ou Coimbra, 25-03-1560 MMHM""}, ""value"": {""comment"": ""@wikidata:Q45412 Domingues & O Neil, IV: 2645. MMHM:p.226"", ""original"": ""?""}}"
$ws.Cells.Item(3, 5).Value = "extra_info: {""value"": {""comment"": ""@wikidata:Q45412 Schutte, Monumenta historica japoniae I.,p.1180"", ""original"": ""?""}}"
$ws.Cells.Item(4, 5).Value = "Corrigido: local de entrada extra_info: {""value"": {""comment"": ""@wikidata:Q45412 Em Coimbra no ano de 1582 segundo a Carta Annua de 1623 BA"", ""original"": ""?""}}"
$ws.Cells.Item(5, 5).Value = "extra_info: {""value"": {""comment"": ""@wikidata:Q45412 MMHM:p.94 (ARSI Japsin 35 13)"", ""original"": ""?""}}"
$ws.Cells.Item(6, 5).Value = "extra_info: {""value"": {""comment"": ""@wikidata:Q45412 Segundo Dominguez & O Neil, V.III,p.2961. Franco Imagem...Coimbra,II,4,c.36,p.575"", ""original"": ""?""}}"
$ws.Cells.Item(7, 5).Value = "extra_info: {""value"": {""comment"": ""@wikidata:Q45412 Segungo Dominguez, J. M., & O\u2019Neill, C. (2001) II, 1113"", ""original"": ""?""}}"
$ws.Cells.Item(8, 5).Value = "extra_info: {""value"": {""comment"": ""@wikidata:Q45412 Dehergne n\u00e3o especifica local, Brockey estudos Filosofia e Teologia em Coimbra"", ""original"": ""?""}}"
$ws.Cells.Item(9, 5).Value = "extra_info: {""date"": {""comment"": ""ou 16090502""}, ""value"": {""comment"": ""@wikidata:Q45412 (Franco, 1719, t.II, p. 612)"", ""original"": ""?""}}"
$ws.Cells.Item(10, 5).Value = "extra_info: {""date"": {""comment"": ""MMHM:p.8""}, ""value"": {""comment"": ""[Adicionado a partir de Franco, Imagem...Coimbra, II, 522] @wikidata:Q45412"", ""original"": ""?""}}"
$ws.Cells.Item(11, 5).Value = "extra_info: {""value"": {""comment"": ""@wikidata:Q45412 \""\""\""Franco, Imagem...Coimbra, v.2 p.616; Barbosa Machado, v.3\""\""\"""", ""original"": ""?""}}"
$ws.Cells.Item(12, 5).Value = "extra_info: {""value"": {""comment"": ""@wikidata:Q45412 MMHM:p.205 (Sebastian da Maia, da Maya, d'Amaya)"", ""original"": ""?""}}"
$ws.Cells.Item(13, 5).Value = "extra_info: {""value"": {""comment"": ""@wikidata:Q45412 Segundo Louis Buglio 1688"", ""original"": ""?""}}"
$ws.Cells.Item(14, 5).Value = "extra_info: {""value"": {""comment"": ""@wikidata:Q45412 MMHM:p.203 (Matias d'Amaia), Schutte, Monumenta historica japoniae I.1234"", ""original"": ""?""}}"
$ws.Cells.Item(15, 5).Value = "extra_info: {""value"": {""comment"": ""@wikidata:Q45412 Ver Brockey, p.230 cit. Carta de A.P. a Viteleschi, de Coimbra, 26 Agosto 1640"", ""original"": ""?""}}"
$ws.Cells.Item(16, 5).Value = "[Brockey diz ""educated at the College of Coimbra"" cap.6 nota 95, citando a carta ânua de 1656, Ajuda, 49-V.14:62v] extra_info: {""value"": {""comment"": ""@wikidata:Q45412 Brockey, cap.6 n.95"", ""original"": ""?""}}"
